# Rename target variable headers in row 1 of Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "exportaciones"
$ws.Range("C1").Value = "inversiones"
$ws.Range("D1").Value = "gasto_publico"
$ws.Range("E1").Value = "consumo_privado"
$ws.Range("F1").Value = "importaciones"
$ws.Range("G1").Value = "variacion"

$wb.Save()
